# Add Category Name, Supplier Name, Reorder Point, and Overstock Threshold
# columns to the Products sheet (shifting the existing Supplier ID column
# from H to I), and populate sample data for every existing product row.

$wb = $excel.ActiveWorkbook
$catWs = $wb.Worksheets.Item("Categories")
$supWs = $wb.Worksheets.Item("Suppliers")
$prodWs = $wb.Worksheets.Item("Products")

# Safety guard: if this migration already ran (header H1 already updated),
# don't run it again and corrupt the data.
if ($prodWs.Cells.Item(1, 8).Value2 -eq "Category Name") {
    return
}

# Build a Category ID -> Category Name lookup from the Categories sheet.
$catMap = @{}
$r = 2
while ($catWs.Cells.Item($r, 1).Value2 -ne $null -and $catWs.Cells.Item($r, 1).Value2 -ne "") {
    $catId = $catWs.Cells.Item($r, 1).Value2
    $catName = $catWs.Cells.Item($r, 2).Value2
    $catMap[$catId] = $catName
    $r = $r + 1
}

# Build a Supplier ID -> Supplier Name lookup from the Suppliers sheet.
$supMap = @{}
$r = 2
while ($supWs.Cells.Item($r, 1).Value2 -ne $null -and $supWs.Cells.Item($r, 1).Value2 -ne "") {
    $supId = $supWs.Cells.Item($r, 1).Value2
    $supName = $supWs.Cells.Item($r, 2).Value2
    $supMap[$supId] = $supName
    $r = $r + 1
}

# Find the last used row in the Products sheet (based on column A / ID).
$lastRow = 1
$r = 2
while ($prodWs.Cells.Item($r, 1).Value2 -ne $null -and $prodWs.Cells.Item($r, 1).Value2 -ne "") {
    $lastRow = $r
    $r = $r + 1
}

# New header row: H = Category Name, I = Supplier ID (moved from old H),
# J = Supplier Name, K = Reorder Point, L = Overstock Threshold.
$prodWs.Cells.Item(1, 8).Value = "Category Name"
$prodWs.Cells.Item(1, 9).Value = "Supplier ID"
$prodWs.Cells.Item(1, 10).Value = "Supplier Name"
$prodWs.Cells.Item(1, 11).Value = "Reorder Point"
$prodWs.Cells.Item(1, 12).Value = "Overstock Threshold"

for ($row = 2; $row -le $lastRow; $row++) {
    $categoryId = $prodWs.Cells.Item($row, 7).Value2
    $supplierId = $prodWs.Cells.Item($row, 8).Value2

    $categoryName = $null
    if ($categoryId -ne $null -and $categoryId -ne "") {
        $categoryName = $catMap[$categoryId]
    }

    $hasSupplier = ($supplierId -ne $null -and $supplierId -ne "")
    $supplierName = $null
    if ($hasSupplier) {
        $supplierName = $supMap[$supplierId]
    }

    # H: Category Name
    if ($categoryName -ne $null) {
        $prodWs.Cells.Item($row, 8).Value = $categoryName
    }

    if ($hasSupplier) {
        # I: Supplier ID (moved from the old column H)
        $prodWs.Cells.Item($row, 9).Value = $supplierId
        # J: Supplier Name
        if ($supplierName -ne $null) {
            $prodWs.Cells.Item($row, 10).Value = $supplierName
        }
        # K / L: inventory thresholds for stocked products
        $prodWs.Cells.Item($row, 11).Value = 10
        $prodWs.Cells.Item($row, 12).Value = 100
    } else {
        # Services have no supplier / inventory thresholds
        $prodWs.Cells.Item($row, 11).Value = 0
        $prodWs.Cells.Item($row, 12).Value = 0
    }
}
